$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 112, pushing the existing rows 112-158 down to 113-159.
$ws.Rows(112).Insert()

# Populate the newly inserted row 112 with the new record's data.
$ws.Range("A112").Value = 5
$ws.Range("B112").Value = "Macroferia Regional de Talca"
$ws.Range("C112").Value = "Maule"
$ws.Range("D112").Value = 45009
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = 100112001
$ws.Range("G112").Value = "Berenjena"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 200
$ws.Range("K112").Value = 8000
$ws.Range("L112").Value = 8000
$ws.Range("M112").Value = 8000
$ws.Range("N112").Value = "$/caja 50 unidades"
$ws.Range("O112").Value = "Región del Maule"
$ws.Range("P112").Value = 160
$ws.Range("Q112").Value = 50
$ws.Range("R112").Value = "Hortaliza"
